$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "provincia" column (F) was previously annotated as an SDMX dimension;
# it is now recategorized as an IAEST measure, consistent with the new
# curated dimensions mentioned in the commit message.
$ws.Range("F2").Value = "iaest-measure:provincia"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"
